# Auto-generated edit script applying the cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.407.20"
$ws.Range("E2").Value = "  +7.31%  "

$ws.Range("D3").Value = "2.381.94"
$ws.Range("E3").Value = "  +4.57%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.76"
$ws.Range("E5").Value = "  +9.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "317.62"
$ws.Range("E6").Value = "  +3.12%  "

$ws.Range("E7").Value = "  +2.26%  "

$ws.Range("E8").Value = "  -0.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.627"
$ws.Range("E9").Value = "  +4.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.68"
$ws.Range("E10").Value = "  +11.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0932"
$ws.Range("E11").Value = "  +4.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.70"
$ws.Range("E12").Value = "  +6.52%  "

$ws.Range("E13").Value = "  +1.75%  "

$ws.Range("E14").Value = "  +5.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.84"
$ws.Range("E15").Value = "  +5.11%  "

$ws.Range("D16").Value = "2.742.48"
$ws.Range("E16").Value = "  +4.57%  "

$ws.Range("D17").Value = "2.378.05"
$ws.Range("E17").Value = "  +3.99%  "

$ws.Range("D18").Value = "45.312.19"
$ws.Range("E18").Value = "  +6.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.62"
$ws.Range("E19").Value = "  +5.53%  "

$ws.Range("E20").Value = "  +4.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.42"
$ws.Range("E21").Value = "  +1.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.73"
$ws.Range("E22").Value = "  +3.03%  "

$ws.Range("E23").Value = "  +5.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.32"
$ws.Range("E24").Value = "  +3.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  +10.34%  "

$ws.Range("E26").Value = "  -0.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.30"
$ws.Range("E27").Value = "  +6.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.51"
$ws.Range("E28").Value = "  +7.46%  "

$ws.Range("E29").Value = "  -0.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.34"
$ws.Range("E30").Value = "  +10.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.94"
$ws.Range("E31").Value = "  +4.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0950"
$ws.Range("E32").Value = "  +12.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "170.37"
$ws.Range("E33").Value = "  +3.99%  "

$ws.Range("E34").Value = "  +15.94%  "

$ws.Range("E35").Value = "  +2.92%  "

$ws.Range("E36").Value = "  +8.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.95"
$ws.Range("E37").Value = "  +11.04%  "

$ws.Range("E38").Value = "  +13.59%  "

$ws.Range("E39").Value = "  +5.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.99"
$ws.Range("E40").Value = "  +10.55%  "

$ws.Range("E41").Value = "  +11.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.27"
$ws.Range("E42").Value = "  +6.46%  "

$ws.Range("E43").Value = "  +7.25%  "

$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.39"
$ws.Range("E44").Value = "  +12.37%  "

$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "71.51"
$ws.Range("E45").Value = "  +4.23%  "

$ws.Range("E46").Value = "  -0.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "116.82"
$ws.Range("E47").Value = "  +6.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.78"
$ws.Range("E48").Value = "  +13.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.66"
$ws.Range("E49").Value = "  +20.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.34"
$ws.Range("E50").Value = "  +8.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.224"
$ws.Range("E51").Value = "  +18.64%  "
